# Update automàtic: dades i banners [2026-02-05 08:09]
# Refresh DATA_EXTRACCIO timestamps and the latest meteo.cat readings for
# every station row in the daily summary sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (YT - Alt Àneu - Bonabé) ---
$ws.Range("E2").Value = "2026-02-05 08:07:51"
$ws.Range("H2").Value = "'95%"
$ws.Range("N2").Value = "-3.0 °C 4:08 TU"
$ws.Range("O2").Value = "-2.7 °C"

# --- Row 3 (Z1 - Alt Àneu - Bonaigua) ---
$ws.Range("E3").Value = "2026-02-05 08:07:54"
$ws.Range("G3").Value = "179 cm"
$ws.Range("H3").Value = "'67%"
$ws.Range("I3").Value = "0.0 mm"
$ws.Range("K3").Value = "0.0 MJ/m2"
$ws.Range("L3").Value = "23.0 km/h - 269º 0:18 TU"
$ws.Range("M3").Value = "-1.2 °C 4:23 TU"
$ws.Range("N3").Value = "-3.9 °C 2:10 TU"
$ws.Range("O3").Value = "-2.5 °C"

# --- Row 4 (DN - Anglès) ---
$ws.Range("E4").Value = "2026-02-05 08:07:56"

# --- Row 5 (DJ - Banyoles) ---
$ws.Range("E5").Value = "2026-02-05 08:07:59"

# --- Row 6 (X4 - Barcelona - el Raval) ---
$ws.Range("E6").Value = "2026-02-05 08:08:01"

# --- Row 7 (D5 - Barcelona - Observatori Fabra) ---
$ws.Range("E7").Value = "2026-02-05 08:08:04"
$ws.Range("H7").Value = "'73%"
$ws.Range("I7").Value = "0.0 mm"
$ws.Range("J7").Value = "994.5 hPa"
$ws.Range("K7").Value = "0.0 MJ/m2"
$ws.Range("L7").Value = "36.4 km/h - 260º 3:03 TU"
$ws.Range("M7").Value = "9.9 °C 0:59 TU"
$ws.Range("N7").Value = "8.7 °C 2:33 TU"
$ws.Range("O7").Value = "9.4 °C"

# --- Row 8 (UN - Cassà de la Selva) ---
$ws.Range("E8").Value = "2026-02-05 08:08:06"

# --- Row 9 (MS - Castellar de N'hug - el Clot del Moro) ---
$ws.Range("E9").Value = "2026-02-05 08:08:09"

# --- Row 10 (W1 - Castelló D'empúries) ---
$ws.Range("E10").Value = "2026-02-05 08:08:11"

# --- Row 11 (DP - Das - Aeròdrom) ---
$ws.Range("E11").Value = "2026-02-05 08:08:13"

# --- Row 12 (XL - El Prat de Llobregat) ---
$ws.Range("E12").Value = "2026-02-05 08:08:16"

# --- Row 13 (VZ - Espolla) ---
$ws.Range("E13").Value = "2026-02-05 08:08:18"

# --- Row 14 (Z7 - Espot) ---
$ws.Range("E14").Value = "2026-02-05 08:08:21"

# --- Row 15 (XJ - Girona) ---
$ws.Range("E15").Value = "2026-02-05 08:08:23"

# --- Row 16 (YU - L'esquirol - Cantonigròs) ---
$ws.Range("E16").Value = "2026-02-05 08:08:26"

# --- Row 17 (CD - La Seu D'urgell) ---
$ws.Range("E17").Value = "2026-02-05 08:08:29"

# --- Row 18 (Z2 - La Vall de Boí - Boí) ---
$ws.Range("E18").Value = "2026-02-05 08:08:31"
$ws.Range("G18").Value = "114 cm"
$ws.Range("H18").Value = "'96%"
$ws.Range("I18").Value = "0.0 mm"
$ws.Range("K18").Value = "0.0 MJ/m2"
$ws.Range("L18").Value = "4.3 km/h - 293º 0:08 TU"
$ws.Range("M18").Value = "-4.5 °C 0:20 TU"
$ws.Range("N18").Value = "-4.6 °C 0:27 TU"
$ws.Range("O18").Value = "-4.5 °C"

# --- Row 19 (VK - Lleida - Raimat) ---
$ws.Range("E19").Value = "2026-02-05 08:08:34"

# --- Row 20 (Z3 - Meranges - Malniu) ---
$ws.Range("E20").Value = "2026-02-05 08:08:36"
$ws.Range("G20").Value = "112 cm"
$ws.Range("H20").Value = "'65%"
$ws.Range("I20").Value = "0.0 mm"
$ws.Range("K20").Value = "0.0 MJ/m2"
$ws.Range("L20").Value = "18.7 km/h - 262º 0:48 TU"
$ws.Range("M20").Value = "-1.1 °C 0:50 TU"
$ws.Range("N20").Value = "-3.5 °C 2:13 TU"
$ws.Range("O20").Value = "-2.2 °C"

# --- Row 21 (YB - Olot) ---
$ws.Range("E21").Value = "2026-02-05 08:08:39"

# --- Row 22 (YP - Palafrugell) ---
$ws.Range("E22").Value = "2026-02-05 08:08:41"

# --- Row 23 (J5 - Pantà de Darnius - Boadella) ---
$ws.Range("E23").Value = "2026-02-05 08:08:44"
$ws.Range("J23").Value = "993.9 hPa"
$ws.Range("L23").Value = "15.8 km/h - 33º 4:04 TU"
$ws.Range("N23").Value = "4.3 °C 5:36 TU"
$ws.Range("O23").Value = "4.9 °C"

# --- Row 24 (D6 - Portbou - Coll dels Belitres) ---
$ws.Range("E24").Value = "2026-02-05 08:08:46"

# --- Row 25 (YA - Puigcerdà) ---
$ws.Range("E25").Value = "2026-02-05 08:08:49"

# --- Row 26 (DG - Queralbs - Núria) ---
$ws.Range("E26").Value = "2026-02-05 08:08:51"

# --- Row 27 (D4 - Roses) ---
$ws.Range("E27").Value = "2026-02-05 08:08:54"

# --- Row 28 (CI - Sant Pau de Segúries) ---
$ws.Range("E28").Value = "2026-02-05 08:08:56"

# --- Row 29 (XS - Santa Coloma de Farners) ---
$ws.Range("E29").Value = "2026-02-05 08:08:59"

# --- Row 30 (ZC - Setcases - Ulldeter) ---
$ws.Range("E30").Value = "2026-02-05 08:09:01"

# --- Row 31 (XH - Sort) ---
$ws.Range("E31").Value = "2026-02-05 08:09:04"
$ws.Range("G31").Value = "1 cm"
$ws.Range("H31").Value = "'96%"
$ws.Range("I31").Value = "0.0 mm"
$ws.Range("J31").Value = "997.1 hPa"
$ws.Range("M31").Value = "3.9 °C 0:26 TU"
$ws.Range("N31").Value = "2.6 °C 3:11 TU"
$ws.Range("O31").Value = "3.3 °C"

# --- Row 32 (XE - Tarragona - Complex Educatiu) ---
$ws.Range("E32").Value = "2026-02-05 08:09:07"

# --- Row 33 (UE - Torroella de Montgrí) ---
$ws.Range("E33").Value = "2026-02-05 08:09:09"

# --- Row 34 (XO - Vic) ---
$ws.Range("E34").Value = "2026-02-05 08:09:11"

# --- Row 35 (VS - Vielha e Mijaran - Lac Redon) ---
$ws.Range("E35").Value = "2026-02-05 08:09:14"
$ws.Range("H35").Value = "'95%"
$ws.Range("L35").Value = "0.0 km/h - 0º 6:30 TU"
$ws.Range("O35").Value = "-3.8 °C"

# --- Row 36 (D7 - Vinebre) ---
$ws.Range("E36").Value = "2026-02-05 08:09:17"
$ws.Range("I36").Value = "0.2 mm"
$ws.Range("J36").Value = "995.8 hPa"
$ws.Range("L36").Value = "8.3 km/h - 74º 5:37 TU"
$ws.Range("M36").Value = "5.9 °C 5:59 TU"
$ws.Range("O36").Value = "5.2 °C"
